$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for the handful of numeric-looking price figures
# that Excel would otherwise silently coerce into a Number; the source
# workbook always stores column D/E figures as literal text.
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "61.072.34"
$ws.Cells.Item(2, 5).Value = "  +0.74%  "

$ws.Cells.Item(3, 4).Value = "2.655.22"
$ws.Cells.Item(3, 5).Value = "  +1.00%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).Value = "577.73"
$ws.Cells.Item(5, 5).Value = "  +0.30%  "

$ws.Cells.Item(6, 4).Value = "144.78"
$ws.Cells.Item(6, 5).Value = "  +1.30%  "

$ws.Cells.Item(7, 4).Value = "0.998"
$ws.Cells.Item(7, 5).Value = "  +0.07%  "

$ws.Cells.Item(8, 5).Value = "  +0.00%  "

$ws.Cells.Item(9, 5).Value = "  +1.99%  "

$ws.Cells.Item(10, 5).Value = "  +1.04%  "

$ws.Cells.Item(11, 4).Value = "0.382"
$ws.Cells.Item(11, 5).Value = "  +4.10%  "

$ws.Cells.Item(12, 5).Value = "  +0.96%  "

$ws.Cells.Item(13, 4).Value = "3.125.05"
$ws.Cells.Item(13, 5).Value = "  +1.14%  "

$ws.Cells.Item(14, 4).Value = "26.02"
$ws.Cells.Item(14, 5).Value = "  +11.77%  "

$ws.Cells.Item(15, 4).Value = "61.064.80"
$ws.Cells.Item(15, 5).Value = "  +0.71%  "

$ws.Cells.Item(16, 5).Value = "  +1.17%  "

$ws.Cells.Item(17, 4).Value = "2.664.15"
$ws.Cells.Item(17, 5).Value = "  +1.66%  "

$ws.Cells.Item(18, 4).Value = "11.65"
$ws.Cells.Item(18, 5).Value = "  +3.37%  "

$ws.Cells.Item(19, 5).Value = "  +1.84%  "

$ws.Cells.Item(20, 4).Value = "351.48"
$ws.Cells.Item(20, 5).Value = "  +0.75%  "

$ws.Cells.Item(21, 5).Value = "  +0.69%  "

$ws.Cells.Item(22, 4).Value = "1.00"
$ws.Cells.Item(22, 5).Value = "  +0.23%  "

$ws.Cells.Item(23, 4).Value = "0.530"
$ws.Cells.Item(23, 5).Value = "  +2.22%  "

$ws.Cells.Item(24, 4).Value = "64.17"
$ws.Cells.Item(24, 5).Value = "  +1.45%  "

$ws.Cells.Item(25, 5).Value = "  +0.76%  "

$ws.Cells.Item(26, 4).Value = "0.997"
$ws.Cells.Item(26, 5).Value = "  +0.10%  "

$ws.Cells.Item(27, 4).Value = "8.19"
$ws.Cells.Item(27, 5).Value = "  +5.18%  "

$ws.Cells.Item(28, 4).Value = "1.97"
$ws.Cells.Item(28, 5).Value = "  +7.24%  "

$ws.Cells.Item(29, 4).Value = "0.0₃0816"
$ws.Cells.Item(29, 5).Value = "  +2.61%  "

$ws.Cells.Item(30, 4).Value = "6.90"
$ws.Cells.Item(30, 5).Value = "  +8.06%  "

$ws.Cells.Item(31, 5).Value = "  +0.05%  "

$ws.Cells.Item(32, 4).Value = "166.20"
$ws.Cells.Item(32, 5).Value = "  +2.70%  "

$ws.Cells.Item(33, 4).Value = "19.98"
$ws.Cells.Item(33, 5).Value = "  +2.13%  "

$ws.Cells.Item(34, 2).Value = "NEARProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(34, 4).Value = "4.52"
$ws.Cells.Item(34, 5).Value = "  +6.89%  "

$ws.Cells.Item(35, 2).Value = "Fetch.AI"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(35, 4).Value = "1.06"
$ws.Cells.Item(35, 5).Value = "  +9.49%  "

$ws.Cells.Item(36, 5).Value = "  +7.48%  "

$ws.Cells.Item(37, 5).Value = "  +4.81%  "

$ws.Cells.Item(38, 4).Value = "338.98"
$ws.Cells.Item(38, 5).Value = "  +12.62%  "

$ws.Cells.Item(39, 4).Value = "4.05"
$ws.Cells.Item(39, 5).Value = "  +4.18%  "

$ws.Cells.Item(40, 4).Value = "0.893"
$ws.Cells.Item(40, 5).Value = "  +5.53%  "

$ws.Cells.Item(41, 4).Value = "38.54"
$ws.Cells.Item(41, 5).Value = "  +1.69%  "

$ws.Cells.Item(42, 5).Value = "  +4.73%  "

$ws.Cells.Item(43, 4).Value = "20.40"

$ws.Cells.Item(44, 4).Value = "134.15"
$ws.Cells.Item(44, 5).Value = "  -0.29%  "

$ws.Cells.Item(45, 2).Value = "VeChain"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(45, 4).Value = "0.0249"
$ws.Cells.Item(45, 5).Value = "  +3.11%  "

$ws.Cells.Item(46, 2).Value = "Stellar"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(46, 4).Value = "0.0999"
$ws.Cells.Item(46, 5).Value = "  +1.41%  "

$ws.Cells.Item(47, 5).Value = "  +2.84%  "

$ws.Cells.Item(48, 4).Value = "0.615"
$ws.Cells.Item(48, 5).Value = "  +1.40%  "

$ws.Cells.Item(49, 4).Value = "20.57"
$ws.Cells.Item(49, 5).Value = "  +3.11%  "

$ws.Cells.Item(50, 5).Value = "  +0.01%  "

$ws.Cells.Item(51, 4).Value = "2.106.13"
$ws.Cells.Item(51, 5).Value = "  +3.93%  "
